$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily rows appended to the liquidity injection/drain series table
# (update fetched from the "MV -datos-" source).
$newRows = @(
    @("28-09-2021", 37, 7, 9, -2, 0, -2, -6, -35, -15, 7),
    @("29-09-2021", 37, 7, 9, -2, 0, -2, -6, -34, -15, 7),
    @("30-09-2021", 36, 7, 8, -2, 0, -2, -6, -35, -13, 7),
    @("01-10-2021", 36, 7, 8, -2, 0, -2, -5, -34, -15, 7),
    @("04-10-2021", 36, 7, 8, -2, 0, -2, -6, -34, -14, 7)
)

$startRow = 191
for ($i = 0; $i -lt $newRows.Length; $i++) {
    $row = $startRow + $i
    $values = $newRows[$i]

    # Column A holds a date string such as "01-10-2021". Assigning it
    # straight to .Value makes Excel's automatic data-type detection treat
    # some of these (day <= 12) as real dates, which would introduce a
    # numeric serial value plus a new number-format style. Build the text
    # with a formula first (never auto-converted) and then convert that
    # formula result to a plain static value, which keeps it as ordinary
    # shared-string text with no style change, matching the rest of the
    # column.
    $dateCell = $ws.Cells.Item($row, 1)
    $dateCell.Formula = '="' + $values[0] + '"'
    $dateCell.Copy()
    $dateCell.PasteSpecial(-4163)
    $excel.CutCopyMode = $false

    for ($c = 1; $c -lt $values.Length; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $values[$c]
    }
}
